$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "Area responsable" text: Secretaria Administrativa -> Departamento
# de Infraestructura, now left-aligned (it picks up the left-align
# formatting that the rest of the row already uses).
# ------------------------------------------------------------------
$ws.Range("H8").Value = "Departamento de Infraestructura (UPP) "
$ws.Range("H8").HorizontalAlignment = -4131   # xlLeft

# ------------------------------------------------------------------
# The empty cells D8:G8 lose their explicit left-alignment (back to
# the default/general alignment).
# ------------------------------------------------------------------
$ws.Range("D8:G8").HorizontalAlignment = 1    # xlGeneral

# ------------------------------------------------------------------
# Roll the reporting period forward: Q4 2022 -> Q2 2023, and bump the
# validation / update dates accordingly.
# ------------------------------------------------------------------
$ws.Range("A8").Value = 2023
$ws.Range("B8").Value = 44927   # 2023-01-01
$ws.Range("C8").Value = 45107   # 2023-06-30
$ws.Range("I8").Value = 45117   # 2023-07-10
$ws.Range("J8").Value = 45117   # 2023-07-10

# Row 8 is a little shorter now that the note text is shorter.
$ws.Rows.Item(8).RowHeight = 45

# Column K (Nota) is a bit wider.
$ws.Columns.Item(11).ColumnWidth = 46.33

# Leave the cursor where the author left it.
$ws.Range("B13").Select()
